# "Generate Report for Archive"
#
# Rows 4 and 5 in every sheet ("Overview", "zh-cn", "de-de") describe the
# two files 4d93103c-...md and f72809ea-...md. The handoff report was
# regenerated and the two files swapped places (f72809ea is now the
# "In Translation" row 4, 4d93103c dropped to the "Ready for handoff" row
# 5), carrying their respective handoff file / datetime columns with them.
# Row formatting (styles) and the hyperlink r:id-to-row mapping stay where
# they are; only the cell values and hyperlink display text move.

$wb = $excel.ActiveWorkbook

function Set-RowHyperlinkDisplay {
    param($ws, $row, $col, $text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $row -and $h.Range.Column -eq $col) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (file), B (status), C (status, mirrored)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

Set-RowHyperlinkDisplay $ws 4 1 "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
Set-RowHyperlinkDisplay $ws 5 1 "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn": A (file), B (status), C (handoff file), D (handoff time)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A4").Value = "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-02 09:38:32"

$ws.Range("A5").Value = "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.zh-cn.xlf"
$ws.Range("D5").Value = "2016-03-02 09:39:16"

Set-RowHyperlinkDisplay $ws 4 1 "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
Set-RowHyperlinkDisplay $ws 4 3 "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.zh-cn.xlf"
Set-RowHyperlinkDisplay $ws 5 1 "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"
Set-RowHyperlinkDisplay $ws 5 3 "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de": A (file), B (status), C (handoff file), D (handoff time)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A4").Value = "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.de-de.xlf"
$ws.Range("D4").Value = "2016-03-02 09:38:42"

$ws.Range("A5").Value = "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.de-de.xlf"
$ws.Range("D5").Value = "2016-03-02 09:39:27"

Set-RowHyperlinkDisplay $ws 4 1 "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.md"
Set-RowHyperlinkDisplay $ws 4 3 "f72809ea-75b7-4db7-bb77-eacc4bb3f8f8.a977ad6808515d0368d3fca445a3c61fedb72c86.de-de.xlf"
Set-RowHyperlinkDisplay $ws 5 1 "4d93103c-a919-4d01-b99a-dbe0d34ebbde.md"
Set-RowHyperlinkDisplay $ws 5 3 "4d93103c-a919-4d01-b99a-dbe0d34ebbde.9d059faf4a0865186050a9f0deda2f5cf2c137a6.de-de.xlf"

Write-Host "Swapped rows 4/5 on Overview, zh-cn, de-de sheets."
